$wb = $excel.ActiveWorkbook

# --- Linear sheet ---
$wsLinear = $wb.Worksheets.Item("Linear")
$wsLinear.Range("B2").Value = -0.0009455211689375831
$wsLinear.Range("B3").Value = -0.01803213825069377
$wsLinear.Range("B4").Value = 0.02113417220723546
$wsLinear.Range("B5").Value = "[1.0, 0.2576584376144352, 0.08811695152619492, 0.0690427323863853, 0.060193913843876116, 0.0981746620280375, 0.2461456960629693, 0.3825555131682569, 0.2487411331636779, 0.09383004411509807, 0.06658773806731169, 0.07734839560176601, 0.09685323362464907, 0.2502705962106188, 0.3993908164000267, 0.23925415163710584, 0.05553930245336261, 0.06530015666677885, 0.06877967134932234, 0.08204819366810377]"

# --- NonLinear sheet ---
$wsNonLinear = $wb.Worksheets.Item("NonLinear")
$wsNonLinear.Range("B4").Value = 0.007334889856758483
$wsNonLinear.Range("B5").Value = 0.1375318038842612
$wsNonLinear.Range("B6").Value = 0.01965054992930092
$wsNonLinear.Range("B7").Value = 0.0006037316301283086
$wsNonLinear.Range("B8").Value = -0.1050824222689125
$wsNonLinear.Range("B9").Value = 0.02262789490044169
$wsNonLinear.Range("B10").Value = "[0.9999999999999999, 0.25905991178147053, 0.09812509907699662, 0.08324290732275103, 0.07244824726966428, 0.11198874356540497, 0.2496101800240498, 0.37729446427876384, 0.24963275799572635, 0.10438087512281577, 0.07905535142367118, 0.08996095564053531, 0.10774471430780701, 0.25148386728802724, 0.39424019419621953, 0.24147224985022683, 0.06616872539595632, 0.07562909679191669, 0.0780327733197378, 0.09115362737889963]"
